$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values reused from the existing shared strings (rows already present at row 43
# establish the same A/B/C pattern that needs to be propagated down through row 57).
$colA = "https://www.quora.com/According-to-the-U-S-News-and-World-Report-rankings-CS-at-UIUC-ranks-better-than-at-prestigious-colleges-like-Caltech-Yale-and-Brown-Is-that-actually-true"
$colB = "Yes"
$colC = "According to the U.S. News and World Report rankings, CS at UIUC ranks better than at prestigious colleges like Caltech, Yale and Brown. Is that actually true?"

# Rows 44-57 only had column D populated; fill in the matching A/B/C (URL / Yes /
# Question) values that mirror row 43 so the whole answer block shares the same
# source row, same as every other question block on the sheet.
for ($r = 44; $r -le 57; $r++) {
    $ws.Cells.Item($r, 1).Value = $colA
    $ws.Cells.Item($r, 2).Value = $colB
    $ws.Cells.Item($r, 3).Value = $colC
}

# Row 37 was re-saved with Excel's real max row height (409.5) instead of 409.6.
$ws.Rows.Item(37).RowHeight = 409.5

# The view had scrolled/selected down to D58; it now shows near the top of the
# newly filled-in block with C44 selected.
$excel.ActiveWindow.ScrollRow = 44
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C44").Select()
